$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: "<Name>_old" -> "<Name>_FV2304" (cols A-J), "diff" stays
# put (col K), "<Name>_new" -> "<Name>_FV2310" (cols L-U).
$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the A1:U58 range into an Excel Table ("Table1") with an autofilter.
$tableRange = $ws.Range("A1:U58")
$listObject = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$listObject.Name = "Table1"

# Freeze the header row (split/freeze at row 2, i.e. below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
